# weitere Bauteile zugeordnet (Widerstaende und Dioden)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- D2..D5 (diodes): LCSC part numbers + "hat Manu" flags, D5 value changed ---
$ws.Range("F30").Value2 = "C205443"
$ws.Range("G30").Value2 = "ja"
$ws.Range("F31").Value2 = "C205445"
$ws.Range("G31").Value2 = "ja"
$ws.Range("F32").Value2 = "C150464"
$ws.Range("G32").Value2 = "ja"
$ws.Range("B33").Value2 = "Orange-3V"
$ws.Range("F33").Value2 = "C205445"
$ws.Range("G33").Value2 = "ja"

# --- R1..R12 (resistors): new values/footprints + LCSC part numbers + flags ---
$ws.Range("B47").Value2 = 470
$ws.Range("C47").Value2 = "Resistor_SMD:R_0402_1005Metric_Pad0.72x0.64mm_HandSolder"
$ws.Range("E47").Value2 = "verwende 510"
$ws.Range("F47").Value2 = "C25170"
$ws.Range("G47").Value2 = "ja"

$ws.Range("E48").Value2 = "verwende 10,7k"
$ws.Range("F48").Value2 = "C22857"
$ws.Range("G48").Value2 = "nein"

$ws.Range("F49").Value2 = "C305012"
$ws.Range("G49").Value2 = "nein (nur 20)"

$ws.Range("B50").Value2 = 270
$ws.Range("C50").Value2 = "Resistor_SMD:R_0402_1005Metric_Pad0.72x0.64mm_HandSolder"
$ws.Range("F50").Value2 = "C310256"
$ws.Range("G50").Value2 = "ja"

$ws.Range("B51").Value2 = 270
$ws.Range("C51").Value2 = "Resistor_SMD:R_0402_1005Metric_Pad0.72x0.64mm_HandSolder"
$ws.Range("F51").Value2 = "C310256"
$ws.Range("G51").Value2 = "ja"

$ws.Range("F52").Value2 = "C304599"
$ws.Range("G52").Value2 = "ja"

$ws.Range("F53").Value2 = "C140214"
$ws.Range("G53").Value2 = "ja"

$ws.Range("F54").Value2 = "C25535"
$ws.Range("G54").Value2 = "ja"

$ws.Range("F55").Value2 = "C294637"
$ws.Range("G55").Value2 = "ja"

$ws.Range("B56").Value2 = 0

$ws.Range("E57").ClearContents()
$ws.Range("B57").Value2 = 270
$ws.Range("C57").Value2 = "Resistor_SMD:R_0402_1005Metric_Pad0.72x0.64mm_HandSolder"

# --- U2, U3, U4, Y1, Y2: "hat Manu" flags ---
$ws.Range("G60").Value2 = "ja"
$ws.Range("G61").Value2 = "nein"
$ws.Range("G62").Value2 = "ja"
$ws.Range("G65").Value2 = "ja"
$ws.Range("G66").Value2 = "ja"

# --- column widths (D narrower, new G column) ---
$ws.Columns.Item(4).ColumnWidth = 49.7
$ws.Columns.Item(7).ColumnWidth = 11.31

# --- restore view selection ---
$ws.Range("F32").Select()
